# "Updated cryptos list" data refresh (GitHub Actions scrape), applied cell-by-cell.
#
# Column D ("Price") holds numeric-looking text, e.g. "1.012" or "0.000008997".
# Assigning such a string straight to .Value lets Excel coerce it to a Double,
# which silently mangles the original formatting (e.g. "0.08940" -> 0.0894,
# "27.57" -> 27.57 but "0.000008997" -> 8.997E-06 internally). To keep these as
# literal text — matching the source spreadsheet, which stores them as strings —
# we prefix with a leading apostrophe, exactly like typing  '1.012  into a cell
# in the Excel UI: Excel strips the apostrophe and keeps the text as-is (quote-prefixed).
# Values that already contain multiple dots (e.g. "27.691.99") or other non-numeric
# characters do not need this treatment, since Excel cannot parse them as numbers anyway.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.691.99'
# Row 3
$ws.Range("D3").Value = '1.847.22'
$ws.Range("E3").Value = '  -0.94%  '
# Row 4
$ws.Range("D4").Value = '''1.012'
$ws.Range("E4").Value = '  -2.62%  '
# Row 5
$ws.Range("D5").Value = '''319.73'
$ws.Range("E5").Value = '  -1.37%  '
# Row 6
$ws.Range("D6").Value = '''1.012'
$ws.Range("E6").Value = '  -2.28%  '
# Row 7
$ws.Range("D7").Value = '''0.4314'
$ws.Range("E7").Value = '  -2.54%  '
# Row 8
$ws.Range("D8").Value = '''0.3745'
$ws.Range("E8").Value = '  -1.51%  '
# Row 9
$ws.Range("D9").Value = '''0.07345'
$ws.Range("E9").Value = '  -1.61%  '
# Row 10
$ws.Range("D10").Value = '''0.8799'
$ws.Range("E10").Value = '  -0.58%  '
# Row 11
$ws.Range("D11").Value = '''21.57'
$ws.Range("E11").Value = '  -0.64%  '
# Row 12
$ws.Range("D12").Value = '1.852.73'
$ws.Range("E12").Value = '  -0.85%  '
# Row 13
$ws.Range("D13").Value = '''6.723'
$ws.Range("E13").Value = '  -0.65%  '
# Row 14
$ws.Range("D14").Value = '''5.448'
$ws.Range("E14").Value = '  -2.07%  '
# Row 15
$ws.Range("D15").Value = '''0.07123'
$ws.Range("E15").Value = '  -1.53%  '
# Row 16
$ws.Range("D16").Value = '''87.79'
$ws.Range("E16").Value = '  +4.76%  '
# Row 17
$ws.Range("E17").Value = '  -2.50%  '
# Row 18
$ws.Range("D18").Value = '''0.000008997'
$ws.Range("E18").Value = '  -1.78%  '
# Row 19
$ws.Range("D19").Value = '''1.011'
$ws.Range("E19").Value = '  -2.35%  '
# Row 20
$ws.Range("D20").Value = '''15.47'
$ws.Range("E20").Value = '  -0.64%  '
# Row 21
$ws.Range("D21").Value = '27.707.17'
$ws.Range("E21").Value = '  -0.25%  '
# Row 22
$ws.Range("D22").Value = '''5.239'
$ws.Range("E22").Value = '  -1.46%  '
# Row 23
$ws.Range("D23").Value = '''11.13'
$ws.Range("E23").Value = '  -2.04%  '
# Row 24
$ws.Range("D24").Value = '2.072.02'
$ws.Range("E24").Value = '  -0.81%  '
# Row 25
$ws.Range("D25").Value = '''2.004'
$ws.Range("E25").Value = '  +0.31%  '
# Row 26
$ws.Range("D26").Value = '''155.76'
$ws.Range("E26").Value = '  -2.14%  '
# Row 27
$ws.Range("E27").Value = '  -1.41%  '
# Row 28
$ws.Range("D28").Value = '''2.124'
$ws.Range("E28").Value = '  +7.15%  '
# Row 29
$ws.Range("D29").Value = '''5.379'
$ws.Range("E29").Value = '  +0.91%  '
# Row 30
$ws.Range("D30").Value = '''120.26'
$ws.Range("E30").Value = '  +1.93%  '
# Row 31
$ws.Range("D31").Value = '''0.08940'
$ws.Range("E31").Value = '  -1.39%  '
# Row 32
$ws.Range("D32").Value = '''1.224'
$ws.Range("E32").Value = '  +0.81%  '
# Row 33
$ws.Range("D33").Value = '''0.7791'
$ws.Range("E33").Value = '  -0.01%  '
# Row 34
$ws.Range("D34").Value = '''4.558'
$ws.Range("E34").Value = '  -0.34%  '
# Row 35
$ws.Range("D35").Value = '''2.916'
$ws.Range("E35").Value = '  -5.82%  '
# Row 36
$ws.Range("E36").Value = '  -2.47%  '
# Row 37
$ws.Range("D37").Value = '''1.140'
$ws.Range("E37").Value = '  -0.97%  '
# Row 38
$ws.Range("D38").Value = '''0.05337'
$ws.Range("E38").Value = '  -0.23%  '
# Row 39
$ws.Range("D39").Value = '''0.01974'
$ws.Range("E39").Value = '  -1.01%  '
# Row 40
$ws.Range("D40").Value = '''7.223'
$ws.Range("E40").Value = '  +4.60%  '
# Row 41
$ws.Range("D41").Value = '''2.872'
$ws.Range("E41").Value = '  +0.58%  '
# Row 42
$ws.Range("D42").Value = '''0.5156'
$ws.Range("E42").Value = '  -0.97%  '
# Row 43
$ws.Range("D43").Value = '''0.1680'
$ws.Range("E43").Value = '  -0.77%  '
# Row 44
$ws.Range("D44").Value = '''8.880'
$ws.Range("E44").Value = '  +2.23%  '
# Row 45
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").Value = '''109.08'
$ws.Range("E45").Value = '  -0.70%  '
# Row 46
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '''10.62'
$ws.Range("E46").Value = '  -0.97%  '
# Row 47
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").Value = '''0.4730'
$ws.Range("E47").Value = '  +0.39%  '
# Row 48
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '''0.06500'
$ws.Range("E48").Value = '  +0.40%  '
# Row 49
$ws.Range("E49").Value = '  -1.47%  '
# Row 50
$ws.Range("D50").Value = '''1.013'
$ws.Range("E50").Value = '  -2.48%  '
# Row 51
$ws.Range("D51").Value = '''1.877'
$ws.Range("E51").Value = '  -2.39%  '
